# ---------------------------------------------------------------------------
# LOM3234.xlsx update
#
#  - Row 10 ("Objetivos:"): the value shown in B10/C10 changes from the old
#    "Estudo de Optica Fisica." text to the first professor entry
#    ("519033 - Carlos Yujiro Shigue").
#  - The long Portuguese paragraphs that used to sit under "Programa resumido:"
#    and "Programa:" (rows 16/18) are deleted, together with the matching
#    "HECHT..." bibliography paragraph (old row 24). The professor rows and the
#    remaining labels are reflowed upward to fill rows 13-24.
#  - The sheet shrinks from A1:C27 down to A1:C24.
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while preserving the bold (col A) / wrap-black (col B) /
# wrap-red (col C) formatting already used throughout the sheet, copying it from a
# row that keeps its original formatting (row 3) so no new cell style is created.
function Set-CellText($address, $text, $formatFromAddress) {
    $ws.Range($address).Value = $text
    $ws.Range($formatFromAddress).Copy() | Out-Null
    $ws.Range($address).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

# "Objetivos:" row (row 10) now shows the first professor instead of the old text.
$ws.Range("B10:C10").Value = "519033 - Carlos Yujiro Shigue"

# Clear the old rows 13-27 completely (values + per-cell formatting) so they can be
# rebuilt with the new layout; Clear() (not ClearContents()) also drops now-unused
# cells instead of leaving empty styled cells behind.
$ws.Range("A13:C27").Clear()

# Rebuild rows 13-24 with the new content, formatting and row heights.
# Row 13
Set-CellText "A13" "Programa resumido:" "A3"
Set-CellText "B13" "519033 - Carlos Yujiro Shigue" "B3"
Set-CellText "C13" "519033 - Carlos Yujiro Shigue" "C3"
$ws.Rows.Item(13).RowHeight = 60

# Row 14
Set-CellText "A14" "Short syllabus:" "A3"
$ws.Rows.Item(14).RowHeight = 60

# Row 15
Set-CellText "A15" "Programa:" "A3"
Set-CellText "B15" "1341653 - Maria José Ramos Sandim" "B3"
Set-CellText "C15" "1341653 - Maria José Ramos Sandim" "C3"
$ws.Rows.Item(15).RowHeight = 120

# Row 16
Set-CellText "A16" "Syllabus:" "A3"
$ws.Rows.Item(16).RowHeight = 120

# Row 17
Set-CellText "A17" "Avaliação:" "A3"
$ws.Rows.Item(17).AutoFit()  # revert to default row height, no custom height

# Row 18
Set-CellText "A18" "Método:" "A3"
Set-CellText "B18" "1643715 - Paulo Atsushi Suzuki" "B3"
Set-CellText "C18" "1643715 - Paulo Atsushi Suzuki" "C3"
$ws.Rows.Item(18).RowHeight = 60

# Row 19
Set-CellText "A19" "Critério:" "A3"
Set-CellText "B19" "Aulas expositivas e práticas ministradas em laboratório." "B3"
Set-CellText "C19" "Aulas expositivas e práticas ministradas em laboratório." "C3"
$ws.Rows.Item(19).RowHeight = 60

# Row 20
Set-CellText "A20" "Norma de recuperação:" "A3"
Set-CellText "B20" "Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4" "B3"
Set-CellText "C20" "Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4" "C3"
$ws.Rows.Item(20).RowHeight = 60

# Row 21
Set-CellText "A21" "Bibliografia:" "A3"
Set-CellText "B21" "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação" "B3"
Set-CellText "C21" "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação" "C3"
$ws.Rows.Item(21).RowHeight = 120

# Row 22
Set-CellText "A22" "Requisitos:" "A3"
$ws.Rows.Item(22).AutoFit()  # revert to default row height, no custom height

# Row 23
Set-CellText "B23" "LOB1021 -  Física IV  (Requisito)`n" "B3"
Set-CellText "C23" "LOB1021 -  Física IV  (Requisito)`n" "C3"
$ws.Rows.Item(23).RowHeight = 30

# Row 24
Set-CellText "B24" "LOM3205 -  Eletromagnetismo  (Requisito)`n" "B3"
Set-CellText "C24" "LOM3205 -  Eletromagnetismo  (Requisito)`n" "C3"
$ws.Rows.Item(24).RowHeight = 30

# Rows 25-27 (old "Requisitos:" block) are now redundant in this position; remove them
# entirely so the sheet dimension shrinks from A1:C27 to A1:C24.
$ws.Range("A25:C27").EntireRow.Delete()

$excel.CutCopyMode = $false

